$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.98
$ws.Range("H2").Value = 2.62
$ws.Range("P2").Value = 2.02
$ws.Range("T2").Value = 1.67
$ws.Range("W2").Value = 1.5
$ws.Range("AF2").Value = 21
$ws.Range("AK2").Value = 36
$ws.Range("Y3").Value = 12.5
$ws.Range("Z3").Value = 16.5
$ws.Range("AB3").Value = 18.5
$ws.Range("AC3").Value = 8.800000000000001
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 25
$ws.Range("AG3").Value = 18.5
$ws.Range("AO3").Value = 1000
$ws.Range("G5").Value = 1.43
$ws.Range("H5").Value = 8.199999999999999
$ws.Range("X5").Value = 21
$ws.Range("AK5").Value = 16
$ws.Range("F6").Value = 5.2
$ws.Range("H6").Value = 1.79
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 3.9
$ws.Range("Q6").Value = 2.06
$ws.Range("H7").Value = 14.5
$ws.Range("O7").Value = 1.14
$ws.Range("S7").Value = 2.08
$ws.Range("T7").Value = 1.99
$ws.Range("U7").Value = 1.88
$ws.Range("X7").Value = 100
$ws.Range("AC7").Value = 24
$ws.Range("F8").Value = 1.78
$ws.Range("G8").Value = 1.79
$ws.Range("X8").Value = 1000
$ws.Range("AK8").Value = 21
$ws.Range("AL8").Value = 1000
$ws.Range("H9").Value = 6.6
$ws.Range("K9").Value = 4.9
$ws.Range("Q9").Value = 1.69
$ws.Range("X9").Value = 30
$ws.Range("Z9").Value = 390
$ws.Range("F10").Value = 1.52
$ws.Range("G10").Value = 1.53
$ws.Range("Z10").Value = 820
$ws.Range("AD10").Value = 42
$ws.Range("F11").Value = 3.3
$ws.Range("G11").Value = 3.4
$ws.Range("H11").Value = 2.34
$ws.Range("I11").Value = 2.42
$ws.Range("AB11").Value = 13.5
$ws.Range("AG11").Value = 15.5
$ws.Range("AJ11").Value = 75
$ws.Range("AN11").Value = 44
$ws.Range("G12").Value = 1.59
$ws.Range("H12").Value = 6
$ws.Range("Y12").Value = 1000
$ws.Range("AD12").Value = 36
$ws.Range("AL12").Value = 1000
$ws.Range("I13").Value = 1.79
$ws.Range("G14").Value = 2.44
$ws.Range("H14").Value = 3.35
$ws.Range("P14").Value = 1.69
$ws.Range("Q14").Value = 2.08
$ws.Range("I15").Value = 17.5
$ws.Range("F17").Value = 2.44
$ws.Range("G17").Value = 2.74
$ws.Range("H17").Value = 2.82
$ws.Range("P17").Value = 1.93
$ws.Range("Q17").Value = 1.73
$ws.Range("G22").Value = 7.4
$ws.Range("H22").Value = 1.51
$ws.Range("Q22").Value = 1.6
$ws.Range("AD22").Value = 11
$ws.Range("AE22").Value = 16
$ws.Range("AF22").Value = 65
$ws.Range("AI22").Value = 34
$ws.Range("F24").Value = 1.99
$ws.Range("I24").Value = 3.95
$ws.Range("P24").Value = 2.38
$ws.Range("U24").Value = 2.48
$ws.Range("Z24").Value = 50
$ws.Range("H25").Value = 5.4
$ws.Range("K25").Value = 4
$ws.Range("Q25").Value = 1.92
$ws.Range("Q26").Value = 1.87
$ws.Range("U27").Value = 1.85
$ws.Range("T28").Value = 1.89
$ws.Range("AL28").Value = 1000
$ws.Range("H29").Value = 16
$ws.Range("K29").Value = 8
$ws.Range("P29").Value = 2.58
$ws.Range("Q29").Value = 1.56
$ws.Range("R29").Value = 1.67
$ws.Range("U29").Value = 1.67
$ws.Range("AB29").Value = 9.6
$ws.Range("AF29").Value = 8
$ws.Range("AG29").Value = 13
$ws.Range("AH29").Value = 100
$ws.Range("AJ29").Value = 9.4
$ws.Range("AL29").Value = 140
$ws.Range("AM29").Value = 300
$ws.Range("AN29").Value = 3.95
$ws.Range("F30").Value = 2.28
$ws.Range("G30").Value = 2.36
$ws.Range("H30").Value = 3.3
$ws.Range("J30").Value = 3.3
$ws.Range("F31").Value = 2.28
$ws.Range("K31").Value = 3.7
$ws.Range("F33").Value = 1.38
$ws.Range("G34").Value = 4.8
$ws.Range("H34").Value = 1.88
$ws.Range("I34").Value = 1.98
$ws.Range("P34").Value = 1.94
$ws.Range("H35").Value = 2.96
$ws.Range("J35").Value = 3.5
$ws.Range("F37").Value = 1.73
$ws.Range("F38").Value = 3.3
$ws.Range("G40").Value = 1.97
